$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output ("NOT FOUND: " + $old)
    }
}

# Collapse runs that were only split apart by proofErr (spelling/grammar) markers.
# The visible text is unchanged; replacing with identical text merges the runs.
Replace-Text "Obrigatório o uso do cinturão de segurança em atividades com altura igual ou maior que dois metros, ou quando o risco da atividade assim determinar;" "Obrigatório o uso do cinturão de segurança em atividades com altura igual ou maior que dois metros, ou quando o risco da atividade assim determinar;"

Replace-Text "Todos os funcionários devem obrigatoriamente fazer uso dos EPI’s - Equipamentos de Proteção Individual fornecidos pela empresa para a finalidade a que se destinam, inclusive roupas adequadas ao frio e umidade;" "Todos os funcionários devem obrigatoriamente fazer uso dos EPI’s - Equipamentos de Proteção Individual fornecidos pela empresa para a finalidade a que se destinam, inclusive roupas adequadas ao frio e umidade;"

Replace-Text "Os funcionários devem responsabilizar-se pelo uso correto, guarda e conservação dos EPI’s;" "Os funcionários devem responsabilizar-se pelo uso correto, guarda e conservação dos EPI’s;"

Replace-Text "Comunicar à Segurança do Trabalho ou seu superior quando os EPI’s se tornarem impróprios para uso, solicitando sua substituição;" "Comunicar à Segurança do Trabalho ou seu superior quando os EPI’s se tornarem impróprios para uso, solicitando sua substituição;"

Replace-Text "Só execute serviços ou opere máquinas se estiver devidamente habilitado e autorizado, quando não souber ou tiver dúvidas sobre algum serviço, pergunte ao seu superior antes do início ou durante a realização do mesmo, para prevenir-se contra possíveis acidentes;" "Só execute serviços ou opere máquinas se estiver devidamente habilitado e autorizado, quando não souber ou tiver dúvidas sobre algum serviço, pergunte ao seu superior antes do início ou durante a realização do mesmo, para prevenir-se contra possíveis acidentes;"

Replace-Text "Em caso de eventuais Acidentes do Trabalho, o funcionário deve de imediato comunicar a Segurança do Trabalho, chefia e/ou responsável, para que o mesmo receba os primeiros socorros e seja feita a abertura da Comunicação de Acidentes de Trabalho – CAT; encaminhando à vítima ao ambulatório e/ou o hospital/posto de atendimento mais próximo de acordo com o Plano de Emergência do parque eólico." "Em caso de eventuais Acidentes do Trabalho, o funcionário deve de imediato comunicar a Segurança do Trabalho, chefia e/ou responsável, para que o mesmo receba os primeiros socorros e seja feita a abertura da Comunicação de Acidentes de Trabalho – CAT; encaminhando à vítima ao ambulatório e/ou o hospital/posto de atendimento mais próximo de acordo com o Plano de Emergência do parque eólico."

# Signature block: name, title, and registry number become generic selectable placeholders.
Replace-Text "BRUNA PETRONI CEZARIO" "NOMEHSE"
Replace-Text "Engenheira de Segurança do Trabalho" "Engenheiro(a) de Segurança do Trabalho"
Replace-Text "CREA-RN: 2122993685" "REGISTROHSE"

Write-Output "done"
